# T1048_Contact_CreateNewContactRecordTypeHoulihanEmployee
# The "Users" sheet previously had a single "User" column (header "User",
# value "Jennie Baker"). It is being split into two role-specific columns:
#   Column A -> "AdminUser" users (header renamed, value replaced)
#   Column B -> "HR User" users (new column, reusing the old "Jennie Baker" value)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Column A: rename header and replace the existing admin user's name
$ws.Range("A1").Value = "AdminUser"
$ws.Range("A2").Value = "Indrajeet Singh"

# Column B: new header + the HR user (former sole "Jennie Baker" entry)
$ws.Range("B1").Value = "HR User"
$ws.Range("B2").Value = "Jennie Baker"

# Header row is bold, matching the existing header style
$ws.Range("B1").Font.Bold = $true

# Make "Users" the active sheet / active cell, as it was the last one edited
$ws.Activate()
$ws.Range("D18").Select()
